# Apply the DE_table2_F14_dim10 edit:
#  - Header A1: "Gen" -> "MaxFES"
#  - Column A (rows 2-14): generation counts -> MaxFES fractions
#  - Column AZ (the "Run 50" column) becomes the recalculated "Mean" column
#    (header + the 13 data rows)
#  - Column BA (the old "Mean" column) is removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Header rename: Gen -> MaxFES
$ws.Cells.Item(1, 1).Value = "MaxFES"

# 2) Column A data values (rows 2-14)
$aValues = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $aValues[$i]
}

# 3) AZ column becomes "Mean" (header + recalculated values)
$ws.Cells.Item(1, 52).Value = "Mean"

$meanValues = @(65.79398498, 49.59223693, 3.00319419, 0.23993097, 0.20721726, 0.18477299, 0.1706826, 0.15506882, 0.14683365, 0.14120652, 0.13258661, 0.12738748, 0.12530687)
for ($i = 0; $i -lt $meanValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, 52).Value = $meanValues[$i]
}

# 4) Remove the old "Mean" column (BA) entirely - shifts the used range back to A:AZ
$ws.Range("BA1:BA14").Delete()
